$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell text updates per the crypto price/volume refresh.
# Values that look numeric (e.g. "1.00", "345.00") are written with a
# leading apostrophe so Excel stores them as text (preserving trailing
# zeros / exact formatting) instead of silently coercing to a number,
# then the cell style is reset to Normal so no stray number-format /
# quote-prefix style gets attached to the cell.

$ws.Range('D2').Value = '63.809.70'
$ws.Range('E2').Value = '  -1.03%  '
$ws.Range('D3').Value = '2.643.73'
$ws.Range('E3').Value = '  +0.54%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = "'580.62"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.45%  '
$ws.Range('D6').Value = "'155.78"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -4.10%  '
$ws.Range('D9').Value = '2.641.39'
$ws.Range('E9').Value = '  +0.60%  '
$ws.Range('E10').Value = '  -3.19%  '
$ws.Range('E11').Value = '  +0.43%  '
$ws.Range('E12').Value = '  -1.65%  '
$ws.Range('E13').Value = '  +0.89%  '
$ws.Range('E14').Value = '  +0.07%  '
$ws.Range('D15').Value = '3.120.19'
$ws.Range('E15').Value = '  +0.61%  '
$ws.Range('D16').Value = "'0.0000184"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.20%  '
$ws.Range('D17').Value = '63.818.73'
$ws.Range('D18').Value = '2.646.63'
$ws.Range('E18').Value = '  -0.56%  '
$ws.Range('E19').Value = '  -0.59%  '
$ws.Range('E20').Value = '  +3.94%  '
$ws.Range('E21').Value = '  -2.90%  '
$ws.Range('D22').Value = "'345.00"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.21%  '
$ws.Range('E23').Value = '  +0.38%  '
$ws.Range('D24').Value = "'68.03"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.44%  '
$ws.Range('D25').Value = "'1.89"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +9.25%  '
$ws.Range('E26').Value = '  -2.36%  '
$ws.Range('D27').Value = "'603.66"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +8.36%  '
$ws.Range('E28').Value = '  -0.70%  '
$ws.Range('E29').Value = '  +3.33%  '
$ws.Range('D30').Value = "'8.12"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.82%  '
$ws.Range('E31').Value = '  -0.12%  '
$ws.Range('D32').Value = "'1.00"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.03%  '
$ws.Range('D33').Value = "'2.06"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('E34').Value = '  +1.88%  '
$ws.Range('D35').Value = "'6.62"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('D36').Value = "'5.45"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.32%  '
$ws.Range('E37').Value = '  -1.65%  '
$ws.Range('E38').Value = '  -1.14%  '
$ws.Range('D39').Value = "'0.999"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E40').Value = '  -1.07%  '
$ws.Range('D41').Value = "'151.25"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.90%  '
$ws.Range('E42').Value = '  +5.56%  '
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('D44').Value = "'41.93"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('D45').Value = "'160.40"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.17%  '
$ws.Range('D46').Value = "'24.58"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +8.10%  '
$ws.Range('E47').Value = '  -1.55%  '
$ws.Range('D48').Value = "'0.0588"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.20%  '
$ws.Range('E49').Value = '  +0.15%  '
$ws.Range('D50').Value = "'0.0999"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.89%  '
$ws.Range('E51').Value = '  -0.72%  '
